$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous contents (old A1:D5 data) so nothing stale remains.
$ws.Range("A1:D5").Clear()

$headers = @("列名1", "列名2", "列名3", "列名4", "列名5", "列名6", "列名7")

# First table: header row 1, data rows 2-4 in column E only.
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$ws.Range("E2").Value = 1234
$ws.Range("E3").Value = 1235
$ws.Range("E4").Value = 1236

# Second table: header row 7, data rows 8-11 in column E only.
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(7, $i + 1).Value = $headers[$i]
}
$ws.Range("E8").Value = 1034
$ws.Range("E9").Value = 1035
$ws.Range("E10").Value = 1036
$ws.Range("E11").Value = 1037
